# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet right before the "总计" (total) sheet,
#    populated with the quarter's fund-holding breakdown.
# 2. Update the "总计" sheet with a new leading row summarizing 2022-Q1
#    (shifting the existing quarters down).

$wb = $excel.ActiveWorkbook

# Reuse an existing quarter sheet purely as a formatting template (same bold /
# bordered header style, same "index" column styling) so the new sheet visually
# matches its siblings.
$template = $wb.Worksheets.Item("2021-Q3")

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# IMPORTANT: $total was resolved to the sheet *positioned* right before "总计";
# Worksheets.Add(Before:=$total) inserts the new sheet into that slot, so the
# $total handle now actually tracks the freshly added "2022-Q1" sheet (handles
# here follow tab position, not the original identity). Re-resolve "总计" by
# name now that the sheet order is final.
$total = $wb.Worksheets.Item("总计")

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1

# Match the page setup used by the other quarter sheets (0.75"/1"/0.5" margins,
# vs. the engine's own new-sheet defaults). PageSetup margins are in points
# (72pt = 1in).
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Match the header/index-column formatting used by the other quarter sheets.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

# Force text format on the numeric-looking fields so codes/figures keep their
# exact literal text (leading zeros, fixed decimal places) instead of being
# coerced to numbers.
$q1Data = $q1.Range("B2:G3")
$q1Data.NumberFormat = "@"

$q1.Range("B2").Value = "001110"
$q1.Range("C2").Value = "中欧瑾泉灵活配置混合 - A"
$q1.Range("D2").Value = "7.39"
$q1.Range("E2").Value = "22.25"
$q1.Range("F2").Value = "1.10"
$q1.Range("G2").Value = "0.0813"
$q1.Range("H2").Value = 10

$q1.Range("B3").Value = "001111"
$q1.Range("C3").Value = "中欧瑾泉灵活配置混合 - C"
$q1.Range("D3").Value = "2.20"
$q1.Range("E3").Value = "22.25"
$q1.Range("F3").Value = "1.10"
$q1.Range("G3").Value = "0.0242"
$q1.Range("H3").Value = 10

# Put the data cells back on the plain/default formatting used elsewhere in the
# workbook (drops the forced "@" number-format residue, keeps the literal text).
$template.Range("B2:G3").Copy()
$q1Data.PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. "总计" sheet: add the 2022-Q1 summary row on top, push the rest down
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.11

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# The blank row Insert() creates picks up neighbouring formatting inconsistently;
# normalize it to match the rest of the table (index column styled, data columns
# plain).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)
